# Trade #11 closed at 2026-02-16 22:53:03 - base_strategy UP +0.000%
#
# The "All Trades" and "base_strategy" sheets are running trade logs.
# Both get a new row appended (row 12) describing the newly-opened
# trade #11.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 12

    $ws.Cells.Item($row, 1).Value = 11

    # Columns B ("Date") and C ("Time") hold plain text that happens to
    # look like a date/time ("2026-02-16", "22:53:03"). Pre-set the
    # number format to Text so the COM layer stores the literal string
    # instead of silently converting it to a date/time serial number,
    # then reset the style back to Normal afterwards so the cell ends
    # up with the same (default) styling as every other row.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-16"
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = "22:53:03"
    $ws.Cells.Item($row, 3).Style = "Normal"

    $ws.Cells.Item($row, 4).Value = "base_strategy"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 49.999998
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = ""
    $ws.Cells.Item($row, 17).Value = 0
}
